# BIS-1002: Fixed XLS export tests
# Adds a new "Internal Assignment" column (O) to both SAMPLE_TYPE property
# tables on the sheet (the ANTIBODY block starting at row 4, and the VIRUS
# block starting at row 12), mirroring the other boolean flag columns
# (Mandatory / Show in edit views / Multivalued / Unique).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-HeaderCell($addr) {
    $r = $ws.Range($addr)
    $r.Value = "Internal Assignment"
    $r.Font.Name = "Calibri"
    $r.Font.Size = 12
    $r.Font.Bold = $true
    $r.Font.Color = 0
    $r.NumberFormat = "General"
}

function Set-FlagCell($addr) {
    $r = $ws.Range($addr)
    $r.Value = "'FALSE"
    $r.Font.Name = "Calibri"
    $r.Font.Size = 11
    $r.Font.Bold = $false
    $r.Font.Family = 0
    $r.NumberFormat = "General"
}

# Header row cells (property-table header rows 4 and 12)
Set-HeaderCell "O4"
Set-HeaderCell "O12"

# Data rows for the first block (rows 5-7) and second block (rows 13-15)
Set-FlagCell "O5"
Set-FlagCell "O6"
Set-FlagCell "O7"

Set-FlagCell "O13"
Set-FlagCell "O14"
Set-FlagCell "O15"
